# Rename the single data sheet from "Property1" to "DataNode" (unifying the
# DataNode/DataTable/Entity naming) and leave the cursor parked on the cell
# the author last had selected (D26) when the file was saved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "DataNode"
$ws.Range("D26").Select()
